$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates (Receptor / Edge expression values recomputed with new TPM)
$ws.Range("M2").Value = 1.193633666666667
$ws.Range("N2").Value = 3.580901
$ws.Range("O2").Value = 0.0852504197348203
$ws.Range("P2").Value = 0.08525041973482028
$ws.Range("Q2").Value = 1.571096441076667
$ws.Range("R2").Value = 14.13986796969
$ws.Range("S2").Value = 0.0852504197348203
$ws.Range("T2").Value = 0.08525041973482028

# Row 3 updates
$ws.Range("O3").Value = 0.6175422122064692
$ws.Range("P3").Value = 0.6175422122064691
$ws.Range("S3").Value = 0.6175422122064692
$ws.Range("T3").Value = 0.6175422122064691

# Row 4 updates
$ws.Range("M4").Value = 4.028899666666667
$ws.Range("N4").Value = 12.086699
$ws.Range("O4").Value = 0.2877477380576656
$ws.Range("P4").Value = 0.2877477380576655
$ws.Range("Q4").Value = 5.302958608256667
$ws.Range("R4").Value = 47.72662747431
$ws.Range("S4").Value = 0.2877477380576656
$ws.Range("T4").Value = 0.2877477380576655

# Row 5 updates
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.132449
$ws.Range("N5").Value = 0.397347
$ws.Range("O5").Value = 0.009459630001044888
$ws.Range("P5").Value = 0.009459630001044887
$ws.Range("Q5").Value = 0.17433334727
$ws.Range("R5").Value = 1.56900012543
$ws.Range("S5").Value = 0.009459630001044888
$ws.Range("T5").Value = 0.009459630001044887
